# "Added Help mark in Home Page"
#
# CreateSTP (sheet1): rename the "17016" full/short-name test fixtures to
# "21012", and extend the "On Par" competitor label to "On Par with
# competation".
#
# CreateSTP_Mandatory (sheet3): append a new Ferrari19/Ferrari141 data row
# (row 21) that mirrors the existing mandatory-fields rows above it.
#
# Selections are left wherever the author's last click landed on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: CreateSTP
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("CreateSTP")

$ws1.Range("F2").Value = "FullName121012"
$ws1.Range("G2").Value = "FullName221012"
$ws1.Range("F3").Value = "Short121012"
$ws1.Range("G3").Value = "Short221012"
$ws1.Range("G26").Value = "On Par with competation"

# ---------------------------------------------------------------------
# Sheet 2: PostDetails (no cell content changes, selection only)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("PostDetails")
$ws2.Activate() | Out-Null
$ws2.Range("M18").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: CreateSTP_Mandatory
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CreateSTP_Mandatory")

$ws3.Range("A21").Value = "Ferrari19"
$ws3.Range("B21").Value = "Ferrari141"
$ws3.Range("C21").Value = "This is a valid Mandatory Fields scenario. It is supposed give a toast message ""STP created successfully!"""

$ws3.Activate() | Out-Null
$ws3.Range("D21").Select() | Out-Null

# Leave CreateSTP as the active/selected tab when the workbook is saved
# (it was the only sheet with tabSelected="1" before the edit, and the
# diff does not touch that attribute), while still landing the cursor on
# E7 as its last-used selection.
$ws1.Activate() | Out-Null
$ws1.Range("E7").Select() | Out-Null
